$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'28.636.89"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = "'  +1.07%  "
$ws.Range("E2").ClearFormats()

$ws.Range("D3").Value = "'1.865.73"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = "'  +1.37%  "
$ws.Range("E3").ClearFormats()

$ws.Range("E4").Value = "'  -0.05%  "
$ws.Range("E4").ClearFormats()

$ws.Range("D5").Value = "'326.29"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "'  -2.27%  "
$ws.Range("E5").ClearFormats()

$ws.Range("D6").Value = "'1.005"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "'  -0.01%  "
$ws.Range("E6").ClearFormats()

$ws.Range("D7").Value = "'0.4659"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "'  +0.65%  "
$ws.Range("E7").ClearFormats()

$ws.Range("D8").Value = "'0.3914"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "'  +0.87%  "
$ws.Range("E8").ClearFormats()

$ws.Range("D9").Value = "'0.07895"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "'  +0.32%  "
$ws.Range("E9").ClearFormats()

$ws.Range("D10").Value = "'0.9740"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "'  +0.46%  "
$ws.Range("E10").ClearFormats()

$ws.Range("D11").Value = "'22.27"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "'  +1.43%  "
$ws.Range("E11").ClearFormats()

$ws.Range("D12").Value = "'1.892.19"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "'  +0.65%  "
$ws.Range("E12").ClearFormats()

$ws.Range("D13").Value = "'5.735"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "'  -0.98%  "
$ws.Range("E13").ClearFormats()

$ws.Range("D14").Value = "'6.952"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "'  +0.22%  "
$ws.Range("E14").ClearFormats()

$ws.Range("D15").Value = "'0.06910"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "'  +0.06%  "
$ws.Range("E15").ClearFormats()

$ws.Range("D16").Value = "'88.84"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "'  +1.59%  "
$ws.Range("E16").ClearFormats()

$ws.Range("E17").Value = "'  -0.03%  "
$ws.Range("E17").ClearFormats()

$ws.Range("D18").Value = "'0.00001002"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "'  +0.22%  "
$ws.Range("E18").ClearFormats()

$ws.Range("D19").Value = "'16.91"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "'  -0.34%  "
$ws.Range("E19").ClearFormats()

$ws.Range("D20").Value = "'1.004"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "'  +0.05%  "
$ws.Range("E20").ClearFormats()

$ws.Range("D21").Value = "'28.637.53"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "'  +0.87%  "
$ws.Range("E21").ClearFormats()

$ws.Range("E22").Value = "'  -0.60%  "
$ws.Range("E22").ClearFormats()

$ws.Range("D23").Value = "'11.07"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "'  -0.59%  "
$ws.Range("E23").ClearFormats()

$ws.Range("D24").Value = "'2.125"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "'  -2.30%  "
$ws.Range("E24").ClearFormats()

$ws.Range("D25").Value = "'2.076.86"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "'  -0.53%  "
$ws.Range("E25").ClearFormats()

$ws.Range("D26").Value = "'155.14"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "'  +0.97%  "
$ws.Range("E26").ClearFormats()

$ws.Range("D27").Value = "'19.29"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "'  -0.08%  "
$ws.Range("E27").ClearFormats()

$ws.Range("D28").Value = "'5.795"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "'  -1.89%  "
$ws.Range("E28").ClearFormats()

$ws.Range("D29").Value = "'1.994"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "'  +0.44%  "
$ws.Range("E29").ClearFormats()

$ws.Range("E30").Value = "'  +1.77%  "
$ws.Range("E30").ClearFormats()

$ws.Range("D31").Value = "'0.09358"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "'  -0.02%  "
$ws.Range("E31").ClearFormats()

$ws.Range("D32").Value = "'0.9393"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "'  -1.35%  "
$ws.Range("E32").ClearFormats()

$ws.Range("D33").Value = "'5.326"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "'  -0.09%  "
$ws.Range("E33").ClearFormats()

$ws.Range("D34").Value = "'1.342"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "'  +0.98%  "
$ws.Range("E34").ClearFormats()

$ws.Range("D35").Value = "'3.345"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "'  -3.40%  "
$ws.Range("E35").ClearFormats()

$ws.Range("D36").Value = "'0.05841"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "'  -3.40%  "
$ws.Range("E36").ClearFormats()

$ws.Range("D37").Value = "'0.02115"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "'  -2.86%  "
$ws.Range("E37").ClearFormats()

$ws.Range("D38").Value = "'1.155"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "'  -0.26%  "
$ws.Range("E38").ClearFormats()

$ws.Range("D39").Value = "'7.889"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "'  +3.74%  "
$ws.Range("E39").ClearFormats()

$ws.Range("D40").Value = "'0.5649"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "'  -0.25%  "
$ws.Range("E40").ClearFormats()

$ws.Range("D41").Value = "'9.971"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "'  -0.67%  "
$ws.Range("E41").ClearFormats()

$ws.Range("D42").Value = "'0.1777"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "'  -0.65%  "
$ws.Range("E42").ClearFormats()

$ws.Range("D43").Value = "'0.07363"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "'  +4.35%  "
$ws.Range("E43").ClearFormats()

$ws.Range("D44").Value = "'11.76"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "'  +0.96%  "
$ws.Range("E44").ClearFormats()

$ws.Range("B45").Value = "RenderToken"
$ws.Range("C45").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D45").Value = "'2.192"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "'  -8.41%  "
$ws.Range("E45").ClearFormats()

$ws.Range("B46").Value = "Decentraland"
$ws.Range("C46").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D46").Value = "'0.5324"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "'  -0.20%  "
$ws.Range("E46").ClearFormats()

$ws.Range("D47").Value = "'1.140"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "'  -6.95%  "
$ws.Range("E47").ClearFormats()

$ws.Range("D48").Value = "'1.851"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "'  -0.41%  "
$ws.Range("E48").ClearFormats()

$ws.Range("D49").Value = "'114.13"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "'  +1.22%  "
$ws.Range("E49").ClearFormats()

$ws.Range("B50").Value = "MXToken"
$ws.Range("C50").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D50").Value = "'2.355"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "'  +0.76%  "
$ws.Range("E50").ClearFormats()

$ws.Range("B51").Value = "PaxDollar"
$ws.Range("C51").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D51").Value = "'1.005"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "'  +0.10%  "
$ws.Range("E51").ClearFormats()
